$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Refresh the "panel_query_time" timestamps (column F) on the "data" sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:19:45.394957",
    "2021-10-05 14:19:45.394965",
    "2021-10-05 14:19:45.394968",
    "2021-10-05 14:19:45.394971",
    "2021-10-05 14:19:45.394974",
    "2021-10-05 14:19:45.394977",
    "2021-10-05 14:19:45.394980",
    "2021-10-05 14:19:45.394982",
    "2021-10-05 14:19:45.394985",
    "2021-10-05 14:19:45.394988",
    "2021-10-05 14:19:45.394991",
    "2021-10-05 14:19:45.394993",
    "2021-10-05 14:19:45.394996",
    "2021-10-05 14:19:45.394999",
    "2021-10-05 14:19:45.395002",
    "2021-10-05 14:19:45.395004",
    "2021-10-05 14:19:45.395007",
    "2021-10-05 14:19:45.395010",
    "2021-10-05 14:19:45.395013",
    "2021-10-05 14:19:45.395016",
    "2021-10-05 14:19:45.395018",
    "2021-10-05 14:19:45.395021",
    "2021-10-05 14:19:45.395024",
    "2021-10-05 14:19:45.395026",
    "2021-10-05 14:19:45.395029",
    "2021-10-05 14:19:45.395032",
    "2021-10-05 14:19:45.395035",
    "2021-10-05 14:19:45.395038",
    "2021-10-05 14:19:45.395040",
    "2021-10-05 14:19:45.395043",
    "2021-10-05 14:19:45.395045",
    "2021-10-05 14:19:45.395048",
    "2021-10-05 14:19:45.395051",
    "2021-10-05 14:19:45.395054",
    "2021-10-05 14:19:45.395057"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" tab (placed after "data") describing the panel
#    query that produced this workbook
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Congenital myaesthenic syndrome"
$meta.Range("C2").Value = 232
# "2.38" must stay text (not become the number 2.38) - force text storage via
# a temporary "@" number format, then restore the plain/default formatting.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.38"
$ws.Range("A1").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$meta.Range("E2").Value = "2021-09-02T17:13:10.476969Z"
$meta.Range("F2").Value = "2021-10-05 14:19:45.391295"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/232/?format=json"

# Match the bold / bordered / centered header style already used by the
# "data" sheet's header row (B1:F1) and its A-column index cells, by
# copying the formatting across instead of re-deriving it, so the same
# style definition is reused.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)

$ws.Range("E1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Keep "data" as the active/selected sheet (unchanged view state), since the
# new tab is purely additive.
$ws.Activate()
